$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns.Item(4).Delete()
